$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "a"
$ws.Range("B3").Value = "b"
$ws.Range("C3").Value = "c"
$ws.Range("D3").Value = "d"

$ws.Range("D3").Select()
